$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these price cells keep their exact text representation
# (trailing zeros, decimal grouping, etc.) instead of Excel auto-converting
# them to numbers when assigned.
$textCells = @('D5','D6','D8','D9','D10','D14','D15','D18','D19','D21','D22','D23','D24','D25','D27','D30','D31','D36','D38','D39','D40','D41','D42','D43','D45','D47','D48','D49','D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '33.949.34'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').Value = '1.782.67'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '221.44'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').Value = '0.552'
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = '31.32'
$ws.Range('E8').Value = '  -5.02%  '
$ws.Range('D9').Value = '0.286'
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').Value = '0.0710'
$ws.Range('E10').Value = '  +4.80%  '
$ws.Range('E11').Value = '  -1.61%  '
$ws.Range('D12').Value = '2.037.84'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '1.779.24'
$ws.Range('D14').Value = '10.51'
$ws.Range('E14').Value = '  -4.73%  '
$ws.Range('D15').Value = '0.627'
$ws.Range('E15').Value = '  -1.19%  '
$ws.Range('D16').Value = '33.929.06'
$ws.Range('E16').Value = '  -2.00%  '
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '68.09'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = '244.98'
$ws.Range('E19').Value = '  -3.49%  '
$ws.Range('D20').Value = '0.0₃0780'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').Value = '0.998'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = '10.69'
$ws.Range('E22').Value = '  +2.62%  '
$ws.Range('D23').Value = '4.08'
$ws.Range('E23').Value = '  -3.79%  '
$ws.Range('D24').Value = '2.10'
$ws.Range('E24').Value = '  -1.23%  '
$ws.Range('D25').Value = '157.42'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').Value = '7.00'
$ws.Range('E27').Value = '  -1.42%  '
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = '0.0523'
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('D31').Value = '3.69'
$ws.Range('E31').Value = '  -1.74%  '
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('E33').Value = '  -2.42%  '
$ws.Range('E34').Value = '  -2.72%  '
$ws.Range('D35').Value = '1.401.59'
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('D36').Value = '0.640'
$ws.Range('E36').Value = '  +1.62%  '
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('D38').Value = '0.0186'
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('D39').Value = '0.935'
$ws.Range('E39').Value = '  +4.08%  '
$ws.Range('D40').Value = '79.35'
$ws.Range('E40').Value = '  -4.72%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '2.72'
$ws.Range('E41').Value = '  -3.07%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').Value = '2.34'
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('D43').Value = '2.10'
$ws.Range('E43').Value = '  +1.77%  '
$ws.Range('D45').Value = '0.0492'
$ws.Range('E45').Value = '  -1.82%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.937.93'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '1.03'
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('D48').Value = '105.25'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').Value = '0.996'
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('D50').Value = '11.77'
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('E51').Value = '  -1.17%  '
